$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B; this shifts the existing
# NCTId..results columns (B..I) one slot to the right (C..J) and
# keeps their formatting/values intact.
$ws.Columns("B").Insert()

# New header: string version of the "statut" emoji column.
$ws.Cells.Item(1, 2).Value = "status_label"

# Find the last used data row from column A (the status emoji column).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Fill in the French text label for each row based on the emoji
# stored in column A.
for ($r = 2; $r -le $lastRow; $r++) {
    $status = $ws.Cells.Item($r, 1).Value2
    if ($status -eq "🟥") {
        $label = "rouge"
    } elseif ($status -eq "🟧") {
        $label = "orange"
    } elseif ($status -eq "🟩") {
        $label = "vert"
    } else {
        $label = ""
    }
    $ws.Cells.Item($r, 2).Value = $label
}
